$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161; existing rows 161.. shift down to 162..
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new record's data
$ws.Cells.Item(161, 1).Value = 3
$ws.Cells.Item(161, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44518
$ws.Cells.Item(161, 5).Value = 5
$ws.Cells.Item(161, 6).Value = 100112031
$ws.Cells.Item(161, 7).Value = "Poroto verde"
$ws.Cells.Item(161, 8).Value = "Magnum"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 73
$ws.Cells.Item(161, 11).Value = 45000
$ws.Cells.Item(161, 12).Value = 46000
$ws.Cells.Item(161, 13).Value = 45479
$ws.Cells.Item(161, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(161, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(161, 16).Value = 1819
$ws.Cells.Item(161, 17).Value = 25
$ws.Cells.Item(161, 18).Value = "Hortaliza"
